$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add an extra centered blank line (sz=28) between the title
#    "DA Python Projet 5" and the existing blank line that precedes
#    "Données publiques de l'OpenFoodFacts".
# ---------------------------------------------------------------------------
$blank28 = $d.Paragraphs(2)
$blank28.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 2) Add an extra centered blank line (sz=24) between
#    "Données publiques de l'OpenFoodFacts" and the existing blank line
#    that precedes "I) Lien vers le projet 5 ...".
#    After step 1, the pre-existing blank (sz=24) paragraph moved to index 5.
# ---------------------------------------------------------------------------
$blank24 = $d.Paragraphs(5)
$blank24.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 3) "La deuxième étape ..." paragraph:
#    "remplaceurs et un menu de choix pour le produit remplaceur d'un
#     côté, et" -> "remplaçants et un menu de choix pour le produit
#     remplaceur d'un côté ; et"
#    Done as two separate, small in-place edits (each bracketed by a
#    bold on/off toggle) so the edited spans land in their own runs,
#    exactly as an interactive edit would leave them, and so regular
#    apostrophes are preserved (Find/Replace's replacement text turns
#    straight quotes into curly ones).
# ---------------------------------------------------------------------------
$rCeur = $d.Content
$rCeur.Find.Execute("ceur")
$rCeur.Text = "çant"
$rCeur.Font.Bold = $true
$rCeur.Font.Bold = $false

$rS = $d.Range($rCeur.End, $rCeur.End + 1)
$rS.Font.Bold = $true
$rS.Font.Bold = $false

$rComma = $d.Content
$rComma.Find.Execute("d'un côté,")
$commaOnly = $d.Range($rComma.End - 1, $rComma.End)
$commaOnly.Text = " ;"
$commaOnly.Font.Bold = $true
$commaOnly.Font.Bold = $false

# ---------------------------------------------------------------------------
# 4) "La troisième étape ..." paragraph:
#    "continue" -> "continu" (typo fix) and drop the trailing sentence
#    about the library.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("en continue.", $true, $false, $false, $false, $false, $true, 1, $false, "en continu.", 2)

$d.Content.Find.Execute(" Une bibliothèque fait aussi partie de l'application concernant le formatage des menus.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# 5) "Enfin, ..." paragraph: drop the trailing sentence about the key
#    error, and give the now-isolated blank paragraph that follows it
#    the same run formatting as its paragraph mark.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" Mais après avoir implémenté cette fonctionnalité, un nouveau problème est apparu : une erreur de clé dans le dictionnaire de produits suite à la recherche empêchait l'enregistrement.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$enfinRange = $d.Content
$enfinRange.Find.Execute("Enfin, lors des tests")
$enfinParaIndex = $enfinRange.Paragraphs.First.Index
$blankAfterEnfin = $d.Paragraphs($enfinParaIndex + 1)
$blankAfterEnfin.Range.Font.NameFarEast = "Times New Roman"
$blankAfterEnfin.Range.Font.NameBi = "Times New Roman"
$blankAfterEnfin.Range.Font.Size = 12
$blankAfterEnfin.Range.Font.SizeBi = 12
